$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "68.726.84"
$ws.Range("E2").Value = "  +1.05%  "

# Row 3
$ws.Range("D3").Value = "3.284.97"
$ws.Range("E3").Value = "  +0.28%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.53%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "186.35"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.96%  "

# Row 7
$ws.Range("E7").Value = "  +0.12%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.601"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.63%  "

# Row 9
$ws.Range("E9").Value = "  -0.41%  "

# Row 10
$ws.Range("E10").Value = "  -0.93%  "

# Row 11
$ws.Range("E11").Value = "  +0.73%  "

# Row 12
$ws.Range("D12").Value = "3.856.31"
$ws.Range("E12").Value = "  +0.43%  "

# Row 13
$ws.Range("E13").Value = "  -0.33%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.62%  "

# Row 15
$ws.Range("D15").Value = "68.764.53"
$ws.Range("E15").Value = "  +1.17%  "

# Row 16
$ws.Range("E16").Value = "  +1.12%  "

# Row 17
$ws.Range("D17").Value = "3.276.92"
$ws.Range("E17").Value = "  +0.38%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.89"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.57%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.80%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "397.23"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +5.15%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.19%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.87"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.92%  "

# Row 23
$ws.Range("E23").Value = "  -0.09%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.521"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.17%  "

# Row 25
$ws.Range("E25").Value = "  +0.51%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.190"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.73%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.81"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.55%  "

# Row 28
$ws.Range("E28").Value = "  -0.07%  "

# Row 29
$ws.Range("E29").Value = "  +0.17%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.76"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "23.13"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.91%  "

# Row 32
$ws.Range("B32").Value = "Fetch.AI"
$ws.Range("C32").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.30"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.75%  "

# Row 33
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.17%  "

# Row 35
$ws.Range("E35").Value = "  -0.56%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "163.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.47%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.03"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +9.58%  "

# Row 38
$ws.Range("E38").Value = "  -3.02%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "26.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.70%  "

# Row 40
$ws.Range("E40").Value = "  -0.66%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.32%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.55"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.21%  "

# Row 43
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.69"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.30%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0694"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.54%  "

# Row 45
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.47"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.51%  "

# Row 46
$ws.Range("D46").Value = "2.652.10"
$ws.Range("E46").Value = "  -1.15%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "342.67"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.70%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0284"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.38"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.32%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "31.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.23%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.997"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.60%  "
